# Generate Report for Handback
# Update timestamps / priority that reflect a newer handback generation run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" for the 83c433f8 and
#     9ae2d3e4 rows (rows 4 & 5) moves forward a minute.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-23 02:16:24"
$wsOverview.Range("G5").Value = "2016-08-23 02:16:24"

# --- zh-cn sheet: Priority flips from "ht" to "mt" for the 83c433f8 /
#     9ae2d3e4 rows, and both handoff/handback datetimes advance.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H4").Value = "2016-08-23 02:16:18"
$wsZhCn.Range("H5").Value = "2016-08-23 02:16:18"
$wsZhCn.Range("K4").Value = "2016-08-23 02:16:43"
$wsZhCn.Range("K5").Value = "2016-08-23 02:16:43"

# --- de-de sheet: handoff datetime matches the Overview column, handback
#     datetime advances.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-23 02:16:24"
$wsDeDe.Range("H5").Value = "2016-08-23 02:16:24"
$wsDeDe.Range("K4").Value = "2016-08-23 02:16:50"
$wsDeDe.Range("K5").Value = "2016-08-23 02:16:50"
